$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp shown above the table
$ws.Range('A1').Value = 'Datos actualizados a 15 de Junio de 2020 a las 02:39'

# Full refreshed country dataset (Pais, Casos totales, Nuevos casos, Casos
# activos, Recuperados, (unused), Casos criticos, Muertes), already sorted
# descending by "Casos totales" as the sheet expects.
$rows = @(
  'Estados Unidos|2162054|19830|867817|1176384|0|326|117853',
  'Brasil|867882|17086|437512|386981|0|598|43389',
  'Rusia|528964|8835|280050|241966|0|119|6948',
  'India|333008|11382|169689|153799|0|321|9520',
  'Reino Unido|295889|1514|0|0|0|36|41698',
  'España|291008|323|0|0|0|0|27136',
  'Italia|236989|338|176370|26274|0|44|34345',
  'Peru|229736|4604|115579|107469|0|190|6688',
  'Alemania|187671|248|172200|6601|0|3|8870',
  'Iran|187427|2472|148674|29916|0|107|8837',
  'Turquia|178239|1562|151417|22015|0|15|4807',
  'Chile|174293|6938|143704|27266|0|222|3323',
  'Francia|157220|407|72859|54954|0|9|29407',
  'Mexico|142690|3494|104975|20843|0|424|16872',
  'Pakistan|139230|6825|51735|84863|0|81|2632',
  'Arabia Saudita|127541|4233|84720|41849|0|40|972',
  'Canada|98787|377|60272|30369|0|39|8146',
  'Banglades|87520|3141|18730|67619|0|32|1171',
  'China|83132|57|78369|129|0|0|4634',
  'Catar|79602|1186|56898|22631|0|3|73',
  'Sudafrica|70038|4302|38531|30027|0|57|1480',
  'Belgica|60029|111|16589|33785|0|5|9655',
  'Bielorrusia|53973|732|30103|23562|0|5|308',
  'Suecia|51614|38|0|0|0|0|4874',
  'Colombia|50939|2193|19822|29450|0|75|1667',
  'Paises Bajos|48783|143|0|0|0|2|6059',
  'Ecuador|46751|395|23064|19791|0|22|3896',
  'Egipto|44598|1618|11931|31092|0|91|1575',
  'Emiratos Arabes Unidos|42294|304|27462|14543|0|1|289',
  'Singapur|40604|407|29589|10989|0|0|26',
  'Indonesia|38277|857|14531|21612|0|43|2134',
  'Portugal|36690|227|22669|12504|0|5|1517',
  'Kuwait|35920|454|26759|8865|0|7|296',
  'Argentina|31577|1282|9564|21180|0|18|833',
  'Ucrania|31154|648|14082|16183|0|9|889',
  'Suiza|31117|23|28800|379|0|0|1938',
  'Polonia|29392|375|14226|13919|0|10|1247',
  'Filipinas|25930|538|5954|18888|0|14|1088',
  'Irlanda|25303|8|22698|899|0|1|1706',
  'Afganistan|24766|664|4725|19570|0|20|471',
  'Oman|23481|1404|8454|14923|0|5|104',
  'Republica Dominicana|22962|390|13320|9050|0|15|592',
  'Rumania|21999|320|15719|4870|0|16|1410',
  'Irak|20209|1259|8121|11481|0|58|607',
  'Panama|20059|0|13759|5871|0|0|429',
  'Israel|19055|83|15375|3380|0|0|300',
  'Barein|18227|514|12818|5367|0|5|42',
  'Bolivia|17842|913|2768|14489|0|26|585',
  'Japon|17429|47|15643|861|0|1|925',
  'Austria|17109|31|16059|373|0|0|677',
  'Armenia|16667|663|6214|10184|0|5|269',
  'Nigeria|16085|403|5220|10445|0|13|420',
  'Kazajistan|14496|258|9188|5231|0|4|77',
  'Serbia|12310|59|11511|545|0|1|254',
  'Dinamarca|12193|54|11068|528|0|0|597',
  'Corea del Sur|12085|34|10718|1090|0|0|277',
  'Ghana|11964|846|4258|7652|0|6|54',
  'Moldavia|11740|281|6623|4711|0|8|406',
  'Argelia|10919|109|7606|2546|0|7|767',
  'Chequia|10024|33|7226|2469|0|1|329',
  'Azerbaiyan|9957|387|5583|4255|0|4|119',
  'Guatemala|9491|509|1804|7320|0|16|367',
  'Marruecos|8793|101|7765|816|0|0|212',
  'Camerun|8681|0|4836|3633|0|0|212',
  'Noruega|8631|3|8138|251|0|0|242',
  'Honduras|8455|323|894|7251|0|4|310',
  'Malasia|8453|8|7346|986|0|1|121',
  'Australia|7320|18|6838|380|0|0|102',
  'Sudan|7220|213|2610|4151|0|12|459',
  'Finlandia|7104|17|6200|578|0|1|326',
  'Nepal|5760|425|974|4767|0|1|19',
  'Senegal|5090|94|3344|1686|0|0|60',
  'Costa de Marfil|5084|236|2505|2534|0|0|45',
  'Uzbekistan|5080|114|3943|1118|0|0|19',
  'Tayikistan|5035|64|3409|1576|0|0|50',
  'Consejo Danes para los Refugiados|4778|54|600|4071|0|1|107',
  'Guinea|4532|48|3234|1273|0|0|25',
  'Republica de Yibuti|4465|16|2950|1472|0|2|43',
  'Haiti|4165|224|24|4071|0|6|70',
  'Luxemburgo|4070|7|3929|31|0|0|110',
  'Hungria|4069|5|2482|1025|0|3|562',
  'Republica de Macedonia|4057|162|1710|2159|0|9|188',
  'El Salvador|3720|117|1837|1811|0|0|72',
  'Kenia|3594|137|1253|2238|0|3|103',
  'Gabon|3463|0|1024|2416|0|0|23',
  'Etiopia|3345|179|545|2743|0|2|57',
  'Bulgaria|3290|24|1730|1386|0|2|174',
  'Tailandia|3135|1|2987|90|0|0|58',
  'Grecia|3121|9|1374|1564|0|0|183',
  'Venezuela|2904|0|487|2393|0|0|24',
  'Bosnia y Herzegovina|2893|0|2119|611|0|0|163',
  'Somalia|2618|39|577|1953|0|1|88',
  'Mayotte|2298|16|1790|479|0|1|29',
  'Kirguistan|2285|78|1791|467|0|0|27',
  'Croacia|2252|1|2134|11|0|0|107',
  'Cuba|2248|10|1948|216|0|0|84',
  'Republica de Africa Central|2057|0|363|1687|0|0|7',
  'Maldivas|2035|22|1311|716|0|0|8',
  'Estonia|1973|0|1705|199|0|0|69',
  'Sri Lanka|1889|5|1287|591|0|0|11',
  'Islandia|1810|2|1796|4|0|0|10',
  'Mali|1809|33|1088|617|0|0|104',
  'Mauritania|1783|101|332|1364|0|4|87',
  'Lituania|1768|5|1427|266|0|0|75',
  'Costa Rica|1715|53|752|951|0|0|12',
  'Sudan del Sur|1693|0|49|1617|0|0|27',
  'Eslovaquia|1548|3|1410|110|0|0|28',
  'Albania|1521|57|1044|441|0|0|36',
  'Nueva Zelanda|1504|0|1482|0|0|0|22',
  'Eslovenia|1495|3|1359|27|0|0|109',
  'Nicaragua|1464|0|953|456|0|0|55',
  'Guinea-Bisau|1460|0|153|1292|0|0|15',
  'Libano|1446|4|868|546|0|0|32',
  'Zambia|1358|1|1122|225|0|1|11',
  'Guinea Ecuatorial|1306|0|200|1094|0|0|12',
  'Paraguay|1289|28|650|628|0|0|11',
  'Madagascar|1272|20|367|895|0|0|10',
  'Guayana Francesa|1255|94|534|718|0|1|3',
  'Sierra Leona|1169|37|680|438|0|0|51',
  'Hong Kong|1110|0|1067|39|0|0|4',
  'Letonia|1097|0|845|224|0|0|28',
  'Tunez|1096|2|998|49|0|0|49',
  'Republica de Chipre|980|0|807|155|0|0|18',
  'Niger|980|0|885|29|0|0|66',
  'Jordania|961|8|682|270|0|0|9',
  'Burkina Faso|894|2|799|42|0|0|53',
  'Georgia|864|13|703|147|0|0|14',
  'Principado de Andorra|853|0|781|21|0|0|51',
  'Republica del Chad|850|2|720|57|0|1|73',
  'Uruguay|848|1|788|37|0|0|23',
  'Cabo Verde|750|24|301|443|0|0|6',
  'Yemen|728|23|53|511|0|4|164',
  'Congo|728|0|221|483|0|0|24',
  'Crucero|712|0|651|48|0|0|13',
  'Uganda|696|2|240|456|0|0|0',
  'San Marino|694|0|520|132|0|0|42',
  'Santo Tome y Principe|661|2|177|472|0|0|12',
  'Malta|649|3|603|37|0|0|9',
  'Jamaica|615|1|420|185|0|0|10',
  'Mozambique|583|30|151|429|0|1|3',
  'Ruanda|582|41|332|248|0|0|2',
  'Malaui|547|18|69|472|0|1|6',
  'Togo|530|0|291|226|0|0|13',
  'Tanzania|509|0|183|305|0|0|21',
  'Reunion|495|6|460|34|0|0|1',
  'Estado de Palestina|492|3|415|74|0|0|3',
  'Suazilandia|490|4|249|237|0|1|4',
  'Liberia|458|12|219|207|0|0|32',
  'Taiwan|443|0|431|5|0|0|7',
  'Benin|442|30|228|208|0|0|6',
  'Libia|418|0|62|348|0|0|8',
  'Zimbabue|383|27|54|325|0|0|4',
  'Mauricio|337|0|325|2|0|0|10',
  'Isla de Man|336|0|312|0|0|0|24',
  'Vietnam|334|0|323|11|0|0|0',
  'Montenegro|325|1|315|1|0|0|9',
  'Birmania|261|0|167|88|0|0|6',
  'Surinam|208|12|9|196|0|0|3',
  'Martinica|202|0|98|90|0|0|14',
  'Mongolia|197|0|98|99|0|0|0',
  'Islas Caimanes|187|0|115|71|0|0|1',
  'Islas Feroe|187|0|187|0|0|0|0',
  'Siria|177|7|74|97|0|0|6',
  'Comoras|176|0|114|60|0|0|2',
  'Gibraltar|176|0|173|3|0|0|0',
  'Guadalupe|171|0|157|0|0|0|14',
  'Guyana|159|0|99|48|0|0|12',
  'Bermudas|142|0|127|6|0|0|9',
  'Brunei|141|0|138|1|0|0|2',
  'Angola|140|2|61|73|0|0|6',
  'Camboya|128|2|125|3|0|0|0',
  'Trinidad yTobago|123|6|109|6|0|0|8',
  'Bahamas|103|0|68|24|0|0|11',
  'Aruba|101|0|98|0|0|0|3',
  'Monaco|99|0|93|2|0|0|4',
  'Eritrea|96|31|39|57|0|0|0',
  'Barbados|96|0|83|6|0|0|7',
  'Burundi|85|0|45|39|0|0|1',
  'Liechtenstein|82|0|55|26|0|0|1',
  'San Martin (Parte Holandesa)|77|0|61|1|0|0|15',
  'Butan|66|4|21|45|0|0|0',
  'Botsuana|60|0|24|35|0|0|1',
  'Polinesia Francesa|60|0|60|0|0|0|0',
  'Macao|45|0|45|0|0|0|0',
  'San Martin (Parte Francesa)|42|0|36|3|0|0|3',
  'Puerto Rico|39|0|1|36|0|0|2',
  'Guam|32|0|0|31|0|0|1',
  'Namibia|32|0|17|15|0|0|0',
  'Gambia|28|0|24|3|0|0|1',
  'San Vicente y las Granadinas|27|0|25|2|0|0|0',
  'Antigua y Barbuda|26|0|20|3|0|0|3',
  'Timor Oriental|24|0|24|0|0|0|0',
  'Granada|23|0|22|1|0|0|0',
  'Curazao|22|0|15|6|0|0|1',
  'Nueva Caledonia|21|0|20|1|0|0|0',
  'Belice|20|0|16|2|0|0|2',
  'Santa Lucia|19|0|18|1|0|0|0',
  'Laos|19|0|19|0|0|0|0',
  'Dominica|18|0|16|2|0|0|0',
  'Fiyi|18|0|18|0|0|0|0',
  'Islas Virgenes de los Estados Unidos|17|0|0|17|0|0|0',
  'San Cristobal y Nieves|15|0|15|0|0|0|0',
  'Groenlandia|13|0|13|0|0|0|0',
  'Islas Malvinas|13|0|13|0|0|0|0',
  'Islas Turcas y Caicos|12|0|11|0|0|0|1',
  'Santa Sede|12|0|12|0|0|0|0',
  'Montserrat|11|0|10|0|0|0|1',
  'Seychelles|11|0|11|0|0|0|0',
  'Sahara Occidental|9|0|8|0|0|0|1',
  'Papua Nueva Guinea|8|0|8|0|0|0|0',
  'Islas Virgenes Britanicas|8|0|7|0|0|0|1',
  'Bonaire, San Eustaquio y Saba|7|0|7|0|0|0|0',
  'San Bartolome|6|0|6|0|0|0|0',
  'Lesoto|4|0|2|2|0|0|0',
  'Anguila|3|0|3|0|0|0|0',
  'San Pedro y Miquelon|1|0|1|0|0|0|0'
)

$data = New-Object 'object[,]' $rows.Count,8
for ($i = 0; $i -lt $rows.Count; $i++) {
    $parts = $rows[$i].Split('|')
    $data[$i,0] = $parts[0]
    for ($j = 1; $j -lt 8; $j++) {
        $data[$i,$j] = [int]$parts[$j]
    }
}
$ws.Range('A4:H219').Value = $data
